$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two names lose their accent marks (Ignacio Fernández Fernández ->
# Ignacio Fernandez Fernandez, Naucé López González -> Nauce Lopez Gonzalez).
# Re-assert the full 4x4 table of values so the shared string table is
# rebuilt cleanly with the new (unaccented) strings.
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "NIF"
$ws.Range("C1").Value = "Email"
$ws.Range("D1").Value = "Colegio"

$ws.Range("A2").Value = "Ignacio Fernandez Fernandez"
$ws.Range("B2").Value = "56378435A"
$ws.Range("C2").Value = "ignacio@uniovi.es"
$ws.Range("D2").Value = 350

$ws.Range("A3").Value = "Nauce Lopez Gonzalez"
$ws.Range("B3").Value = "53678541Z"
$ws.Range("C3").Value = "nauce@uniovi.es"
$ws.Range("D3").Value = 440

$ws.Range("A4").Value = "Jorge Riopedre Vega"
$ws.Range("B4").Value = "48976526C"
$ws.Range("C4").Value = "jorge@uniovi.es"
$ws.Range("D4").Value = 220

# Move the active selection from C2 to A3, matching the saved cursor
# position recorded in the workbook.
$ws.Range("A3").Select()
